# RippleTemplate_Interleaved.xlsx update
#   - add a new "Assay" worksheet (instructions/settings table) at the end
#     of the tab order
#   - tidy up the "Compounds" header row formatting (drop the stray style)
#   - move the active selection to the "Patterns" sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "Assay" sheet after the last existing sheet (Barcodes)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAssay = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsAssay.Name = "Assay"

$wsAssay.Range("A1").Value = "Setting"
$wsAssay.Range("B1").Value = "Value"

$wsAssay.Range("A2").Value = "DMSO Tolerance"
$wsAssay.Range("B2").Value = 0.005

$wsAssay.Range("A3").Value = "Well Volume (µL)"
$wsAssay.Range("B3").Value = 25

$wsAssay.Range("A4").Value = "Backfill (µL)"
$wsAssay.Range("B4").Value = 10

$wsAssay.Range("A5").Value = "Allowed Error"
$wsAssay.Range("B5").Value = 0.1

$wsAssay.Range("A6").Value = "Destination Replicates"
$wsAssay.Range("B6").Value = 1

$wsAssay.Range("A7").Value = "Use Intermediate Plates"
$wsAssay.Range("B7").Value = 1

$wsAssay.Range("A8").Value = "DMSO Normalization"
$wsAssay.Range("B8").Value = 1

# leave the cursor where the author left it on this sheet
$wsAssay.Range("H13").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. Compounds header row (A1:F1) loses its explicit style -> back to
#    the workbook default formatting
# ---------------------------------------------------------------------
$wsCompounds = $wb.Worksheets.Item("Compounds")
$wsCompounds.Range("A1:F1").Style = "Normal"

# ---------------------------------------------------------------------
# 3. Patterns becomes the active tab / selected cell
# ---------------------------------------------------------------------
$wsPatterns = $wb.Worksheets.Item("Patterns")
$wsPatterns.Activate() | Out-Null
$wsPatterns.Range("E10").Select() | Out-Null
